$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Val
    )
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "254.57"
Set-TextValue "E2" "3.55%"
Set-TextValue "D3" "28.12"
Set-TextValue "E3" "-5.70%"
Set-TextValue "D4" "5.288"
Set-TextValue "E4" "2.57%"
Set-TextValue "D5" "0.05849"
Set-TextValue "E5" "1.36%"
Set-TextValue "D6" "6.696"
Set-TextValue "E6" "0.63%"
Set-TextValue "D7" "0.8707"
Set-TextValue "E7" "2.34%"
Set-TextValue "D8" "0.9236"
Set-TextValue "E8" "8.07%"
Set-TextValue "D9" "0.1413"
Set-TextValue "E9" "2.71%"
Set-TextValue "B10" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C10" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D10" "0.03475"
Set-TextValue "E10" "2.93%"
Set-TextValue "B11" "MandalaExchangeToken"
Set-TextValue "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.07110"
Set-TextValue "E11" "0.31%"
Set-TextValue "B12" "BitrueCoin"
Set-TextValue "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D12" "0.03179"
Set-TextValue "E12" "-2.51%"
Set-TextValue "B13" "BitMartToken"
Set-TextValue "C13" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D13" "0.09226"
Set-TextValue "E13" "-1.51%"
Set-TextValue "B14" "BitForexToken"
Set-TextValue "C14" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D14" "0.001553"
Set-TextValue "E14" "1.09%"
Set-TextValue "B15" "One"
Set-TextValue "C15" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D15" "0.0006069"
Set-TextValue "E15" "1.45%"
Set-TextValue "B16" "TigerCash"
Set-TextValue "C16" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D16" "0.005836"
Set-TextValue "E16" "-4.30%"
Set-TextValue "B17" "LEO"
Set-TextValue "C17" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D17" "3.498"
Set-TextValue "E17" "-0.28%"
Set-TextValue "B18" "GateToken"
Set-TextValue "C18" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D18" "3.232"
Set-TextValue "E18" "-0.09%"
Set-TextValue "B19" "BTSEToken"
Set-TextValue "C19" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D19" "2.222"
Set-TextValue "E19" "-0.16%"
Set-TextValue "B20" "BitpandaEcosystemToken"
Set-TextValue "C20" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue "D20" "0.3182"
Set-TextValue "E20" "0.74%"
Set-TextValue "E21" "1.45%"
Set-TextValue "D22" "3.524"
Set-TextValue "E22" "0.58%"
Set-TextValue "D23" "0.04171"
Set-TextValue "E23" "0.91%"
Set-TextValue "D24" "0.1379"
Set-TextValue "E24" "-2.25%"
Set-TextValue "B25" "BitKan"
Set-TextValue "C25" "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D25" "0.001234"
Set-TextValue "E25" "0.53%"
Set-TextValue "B26" "HotbitToken"
Set-TextValue "C26" "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D26" "0.005027"
Set-TextValue "E26" "21.49%"
Set-TextValue "D27" "0.0001201"
Set-TextValue "E27" "-0.02%"
Set-TextValue "D28" "0.0001938"
Set-TextValue "E28" "33.68%"
Set-TextValue "E40" "1.80%"
Set-TextValue "D41" "0.1102"
Set-TextValue "E41" "2.98%"
Set-TextValue "D42" "0.003809"
Set-TextValue "E42" "-33.22%"
Set-TextValue "D43" "0.002339"
Set-TextValue "E43" "1.62%"
Set-TextValue "D44" "0.01012"
Set-TextValue "E44" "19.23%"
Set-TextValue "D45" "0.00005219"
Set-TextValue "E45" "-3.93%"
Set-TextValue "E46" "-0.08%"
Set-TextValue "D47" "0.08754"
Set-TextValue "E47" "23.24%"
Set-TextValue "E48" "-1.81%"
Set-TextValue "D49" "0.00002100"
Set-TextValue "E49" "-0.08%"
Set-TextValue "D50" "0.0002000"
Set-TextValue "E50" "-0.08%"
